$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.674.68'
$ws.Range('E2').Value = '  -0.19%  '
$ws.Range('D3').Value = '3.678.40'
$ws.Range('E3').Value = '  -0.76%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '649.99'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -4.20%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '160.73'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('E9').Value = '  -2.53%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.14'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.19%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.441'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.57%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000231'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.44%  '
$ws.Range('D13').Value = '4.295.86'
$ws.Range('E13').Value = '  -0.85%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '32.62'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.18%  '
$ws.Range('D15').Value = '3.673.39'
$ws.Range('E15').Value = '  -0.64%  '
$ws.Range('D16').Value = '69.685.63'
$ws.Range('E16').Value = '  -0.14%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.118'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.59%  '
$ws.Range('E18').Value = '  +0.10%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '15.89'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.60%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.33'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +4.86%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '470.04'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.76%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.654'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.14%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '79.74'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.02%  '
$ws.Range('D24').Value = '3.824.43'
$ws.Range('E24').Value = '  -0.79%  '
$ws.Range('E25').Value = '  -0.09%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000126'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.26%  '
$ws.Range('E27').Value = '  +0.72%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.76'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -4.57%  '
$ws.Range('E29').Value = '  -2.61%  '
$ws.Range('E30').Value = '  -3.79%  '
$ws.Range('E31').Value = '  +0.05%  '
$ws.Range('E32').Value = '  -2.30%  '
$ws.Range('B33').Value = 'Kaspa'
$ws.Range('C33').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.165'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.42%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '26.71'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.00%  '
$ws.Range('B35').Value = 'NEARProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.43'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.29%  '
$ws.Range('D36').Value = '3.671.28'
$ws.Range('E36').Value = '  -0.66%  '
$ws.Range('E37').Value = '  -2.88%  '
$ws.Range('E38').Value = '  -0.06%  '
$ws.Range('E39').Value = '  -5.54%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '177.90'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +5.80%  '
$ws.Range('E41').Value = '  -0.05%  '
$ws.Range('B42').Value = 'Hedera'
$ws.Range('C42').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0894'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.74%  '
$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.17'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.80%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.930'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.80%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '46.79'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.48%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '29.04'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.88%  '
$ws.Range('E47').Value = '  -1.18%  '
$ws.Range('E48').Value = '  -5.32%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.84'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.19%  '
$ws.Range('E50').Value = '  -5.21%  '
$ws.Range('E51').Value = '  -5.93%  '
